{"js": "// The edit splits a single run of text into two runs that carry identical\n// run formatting (color/theme color, size) \u2014 the visible text is unchanged,\n// but \"Apparence d'un ...\" becomes \"Avoir l'a\" + \"pparence d'un ...\".\n// This happens in two places: \"\u00e9chiquier\" and \"plateau d'\u00e9chec\".\n\nfunction splitOoxml(firstPart, secondPart) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:rPr><w:color w:val=\"156082\" w:themeColor=\"accent1\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\">' +\n    firstPart +\n    '</w:t></w:r><w:r><w:rPr><w:color w:val=\"156082\" w:themeColor=\"accent1\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\">' +\n    secondPart +\n    '</w:t></w:r></w:p></w:body></w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst targets = [\n  { whole: \"Apparence d\\u2019un \u00e9chiquier\", first: \"Avoir l\\u2019a\", second: \"pparence d\\u2019un \u00e9chiquier\" },\n  { whole: \"Apparence d\\u2019un plateau d\\u2019\u00e9chec\", first: \"Avoir l\\u2019a\", second: \"pparence d\\u2019un plateau d\\u2019\u00e9chec\" },\n];\n\nfor (const t of targets) {\n  const results = context.document.body.search(t.whole, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + t.whole);\n  }\n\n  for (const range of results.items) {\n    range.insertOoxml(splitOoxml(t.first, t.second), \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The commit splits a single text run into two runs that share the exact\n# same run formatting (color/theme color, size) - only the run boundary\n# changes, turning:\n#   \"Apparence d'un <X>\"\n# into two adjacent runs:\n#   \"Avoir l'a\"  +  \"pparence d'un <X>\"\n# This happens twice in the document (\"...\u00e9chiquier\" and \"...plateau d'\u00e9chec\").\n\n$d = $word.ActiveDocument\n\n$targets = @(\"Apparence d\u2019un \u00e9chiquier\", \"Apparence d\u2019un plateau d\u2019\u00e9chec\")\n\nforeach ($targetText in $targets) {\n    $searchRange = $d.Content\n    $searchRange.Find.ClearFormatting()\n    $searchRange.Find.Text = $targetText\n    $searchRange.Find.MatchCase = $true\n    $searchRange.Find.MatchWholeWord = $false\n    $searchRange.Find.MatchWildcards = $false\n    $found = $searchRange.Find.Execute()\n\n    if ($found) {\n        # Only touch the leading capital \"A\" - replacing it with \"Avoir l'a\"\n        # turns it into the lower-case \"a\" that starts \"apparence...\". The\n        # untouched remainder (\"pparence d'un ...\") keeps its original,\n        # unmodified run formatting automatically.\n        $leadingA = $d.Range($searchRange.Start, $searchRange.Start + 1)\n        $leadingA.Text = \"Avoir l\u2019a\"\n\n        # Toggling Bold on/off forces the edited text to live in its own run\n        # instead of being silently re-merged with the untouched remainder,\n        # while leaving the final formatting identical to the original run.\n        $leadingA.Bold = $true\n        $leadingA.Bold = $false\n    }\n}\n"}
